# Apply "Atualizacao de bases das ligas" update for Argentina Nacional B
# - Row with id 687 (Quilmes vs Talleres Remedios, B=7702210) data is replaced by a resync
#   that reshuffles several same-kickoff-time matches and drops the trailing scheduled
#   fixture (old last row, B=7845890, CA Atlanta vs CA Aldosivi) entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The very last data row (old row 740, id 738: CA Atlanta vs CA Aldosivi, B=7845890)
# is removed outright, shrinking the used range from AC740 to AC739.
$ws.Rows.Item(740).Delete()

# Row 689
$ws.Range("B689").Value = 7702211
$ws.Range("F689").Value = "CA San Miguel"
$ws.Range("G689").Value = "CA Estudiantes Caseros"
$ws.Range("K689").Value = 2.5
$ws.Range("L689").Value = 2.9
$ws.Range("M689").Value = 2.75
$ws.Range("N689").Value = 2.8
$ws.Range("O689").Value = 2.875
$ws.Range("P689").Value = 2.9
$ws.Range("Q689").Value = 0
$ws.Range("R689").Value = 1.875
$ws.Range("S689").Value = 1.975
$ws.Range("T689").Value = 1.75
$ws.Range("U689").Value = 1.975
$ws.Range("V689").Value = 1.875
$ws.Range("X689").Value = 1.875
$ws.Range("Z689").Value = 0
$ws.Range("AA689").Value = -0
$ws.Range("AC689").Value = 0.875

# Row 690
$ws.Range("B690").Value = 7698357
$ws.Range("F690").Value = "CA San Telmo"
$ws.Range("G690").Value = "Gimnasia y Tiro"
$ws.Range("H690").Value = 1
$ws.Range("I690").Value = 1
$ws.Range("K690").Value = 2.2
$ws.Range("L690").Value = 2.875
$ws.Range("M690").Value = 3.2
$ws.Range("N690").Value = 2.55
$ws.Range("P690").Value = 3.2
$ws.Range("R690").Value = 2
$ws.Range("S690").Value = 1.8
$ws.Range("U690").Value = 1.775
$ws.Range("V690").Value = 2.025
$ws.Range("AB690").Value = 0.3875
$ws.Range("AC690").Value = -0.5

# Row 691
$ws.Range("B691").Value = 7702210
$ws.Range("F691").Value = "Quilmes"
$ws.Range("G691").Value = "Talleres Remedios"
$ws.Range("H691").Value = 0
$ws.Range("I691").Value = 0
$ws.Range("K691").Value = 1.615
$ws.Range("L691").Value = 3.5
$ws.Range("M691").Value = 5
$ws.Range("N691").Value = 1.7
$ws.Range("O691").Value = 3.4
$ws.Range("P691").Value = 5.75
$ws.Range("Q691").Value = -0.75
$ws.Range("R691").Value = 1.975
$ws.Range("S691").Value = 1.825
$ws.Range("T691").Value = 2
$ws.Range("U691").Value = 1.875
$ws.Range("V691").Value = 1.925
$ws.Range("X691").Value = 2.4
$ws.Range("Z691").Value = -1
$ws.Range("AA691").Value = 0.825
$ws.Range("AB691").Value = -1
$ws.Range("AC691").Value = 0.925

# Row 695
$ws.Range("B695").Value = 7698657
$ws.Range("F695").Value = "CA Atlanta"
$ws.Range("G695").Value = "Colon"
$ws.Range("H695").Value = 0
$ws.Range("J695").Value = "D"
$ws.Range("K695").Value = 2.75
$ws.Range("L695").Value = 3.1
$ws.Range("M695").Value = 2.375
$ws.Range("N695").Value = 2.9
$ws.Range("O695").Value = 3.3
$ws.Range("P695").Value = 2.5
$ws.Range("Q695").Value = 0
$ws.Range("R695").Value = 2.125
$ws.Range("S695").Value = 1.75
$ws.Range("T695").Value = 2
$ws.Range("U695").Value = 1.9
$ws.Range("V695").Value = 1.95
$ws.Range("W695").Value = -1
$ws.Range("X695").Value = 2.3
$ws.Range("Z695").Value = 0
$ws.Range("AA695").Value = -0
$ws.Range("AC695").Value = 0.95

# Row 696
$ws.Range("B696").Value = 7698656
$ws.Range("F696").Value = "CA Guemes"
$ws.Range("G696").Value = "San Martin de San Juan"
$ws.Range("L696").Value = 2.875
$ws.Range("N696").Value = 2.875
$ws.Range("O696").Value = 2.9
$ws.Range("P696").Value = 2.8
$ws.Range("Q696").Value = 0
$ws.Range("R696").Value = 1.925
$ws.Range("S696").Value = 1.875
$ws.Range("W696").Value = 1.875
$ws.Range("Z696").Value = 0.925

# Row 697
$ws.Range("B697").Value = 7698545
$ws.Range("F697").Value = "Gimnasia Jujuy"
$ws.Range("G697").Value = "Arsenal de Sarandi"
$ws.Range("L697").Value = 3
$ws.Range("N697").Value = 2.1
$ws.Range("O697").Value = 3.2
$ws.Range("P697").Value = 3.8
$ws.Range("Q697").Value = -0.25
$ws.Range("R697").Value = 1.775
$ws.Range("S697").Value = 2.025
$ws.Range("W697").Value = 1.1
$ws.Range("Z697").Value = 0.7749999999999999

# Row 698
$ws.Range("B698").Value = 7698662
$ws.Range("F698").Value = "Defensores de Belgrano"
$ws.Range("G698").Value = "Nueva Chicago"
$ws.Range("H698").Value = 1
$ws.Range("J698").Value = "H"
$ws.Range("K698").Value = 2.1
$ws.Range("L698").Value = 2.75
$ws.Range("M698").Value = 3.75
$ws.Range("N698").Value = 2.1
$ws.Range("O698").Value = 3.1
$ws.Range("P698").Value = 4
$ws.Range("Q698").Value = -0.5
$ws.Range("R698").Value = 2.1
$ws.Range("S698").Value = 1.775
$ws.Range("T698").Value = 1.75
$ws.Range("U698").Value = 1.775
$ws.Range("V698").Value = 2.1
$ws.Range("W698").Value = 1.1
$ws.Range("X698").Value = -1
$ws.Range("Z698").Value = 1.1
$ws.Range("AA698").Value = -1
$ws.Range("AC698").Value = 1.1

# Row 702
$ws.Range("B702").Value = 7698654
$ws.Range("F702").Value = "Patronato Parana"
$ws.Range("G702").Value = "All Boys"
$ws.Range("H702").Value = 0
$ws.Range("I702").Value = 0
$ws.Range("J702").Value = "D"
$ws.Range("L702").Value = 3.1
$ws.Range("M702").Value = 3
$ws.Range("N702").Value = 1.95
$ws.Range("O702").Value = 3.3
$ws.Range("P702").Value = 4.2
$ws.Range("R702").Value = 1.95
$ws.Range("S702").Value = 1.85
$ws.Range("U702").Value = 1.875
$ws.Range("V702").Value = 1.925
$ws.Range("W702").Value = -1
$ws.Range("X702").Value = 2.3
$ws.Range("Z702").Value = -1
$ws.Range("AA702").Value = 0.8500000000000001
$ws.Range("AB702").Value = -1
$ws.Range("AC702").Value = 0.925

# Row 703
$ws.Range("B703").Value = 7698655
$ws.Range("F703").Value = "CA Alvarado"
$ws.Range("G703").Value = "Guillermo Brown"
$ws.Range("H703").Value = 2
$ws.Range("I703").Value = 1
$ws.Range("J703").Value = "H"
$ws.Range("L703").Value = 3
$ws.Range("M703").Value = 3.2
$ws.Range("N703").Value = 2.05
$ws.Range("O703").Value = 3.25
$ws.Range("P703").Value = 4
$ws.Range("R703").Value = 2.025
$ws.Range("S703").Value = 1.775
$ws.Range("U703").Value = 1.925
$ws.Range("V703").Value = 1.875
$ws.Range("W703").Value = 1.05
$ws.Range("X703").Value = -1
$ws.Range("Z703").Value = 1.025
$ws.Range("AA703").Value = -1
$ws.Range("AB703").Value = 0.925
$ws.Range("AC703").Value = -1

# Row 726
$ws.Range("B726").Value = 7698355
$ws.Range("F726").Value = "Deportivo Madryn"
$ws.Range("G726").Value = "Gimnasia y Tiro"
$ws.Range("K726").Value = 2.5
$ws.Range("L726").Value = 2.9
$ws.Range("N726").Value = 2.375
$ws.Range("P726").Value = 3.6
$ws.Range("Q726").Value = -0.25
$ws.Range("R726").Value = 1.975
$ws.Range("S726").Value = 1.825
$ws.Range("U726").Value = 1.975
$ws.Range("V726").Value = 1.825
$ws.Range("Z726").Value = -0.5
$ws.Range("AA726").Value = 0.4125
$ws.Range("AC726").Value = 0.825

# Row 727
$ws.Range("B727").Value = 7698685
$ws.Range("F727").Value = "Defensores Unidos"
$ws.Range("G727").Value = "Club Atletico Mitre"
$ws.Range("K727").Value = 2.75
$ws.Range("L727").Value = 2.75
$ws.Range("N727").Value = 2.8
$ws.Range("P727").Value = 2.8
$ws.Range("Q727").Value = 0
$ws.Range("R727").Value = 1.875
$ws.Range("S727").Value = 1.925
$ws.Range("U727").Value = 1.875
$ws.Range("V727").Value = 1.925
$ws.Range("Z727").Value = 0
$ws.Range("AA727").Value = -0
$ws.Range("AC727").Value = 0.925

# Row 729
$ws.Range("B729").Value = 7864250
$ws.Range("E729").Value = 45347.70833333334
$ws.Range("F729").Value = "Gimnasia Jujuy"
$ws.Range("G729").Value = "Guillermo Brown"
$ws.Range("K729").Value = 1.727
$ws.Range("L729").Value = 3.3
$ws.Range("M729").Value = 4.333
$ws.Range("N729").Value = 1.8
$ws.Range("O729").Value = 3.3
$ws.Range("P729").Value = 5
$ws.Range("Q729").Value = -0.5
$ws.Range("R729").Value = 1.825
$ws.Range("S729").Value = 2.025
$ws.Range("U729").Value = 1.85
$ws.Range("V729").Value = 2

# Row 730
$ws.Range("B730").Value = 7698687
$ws.Range("F730").Value = "CA San Telmo"
$ws.Range("G730").Value = "Nueva Chicago"
$ws.Range("K730").Value = 2.625
$ws.Range("L730").Value = 2.875
$ws.Range("M730").Value = 2.75
$ws.Range("N730").Value = 2.625
$ws.Range("O730").Value = 2.9
$ws.Range("P730").Value = 3
$ws.Range("Q730").Value = 0
$ws.Range("R730").Value = 1.8
$ws.Range("S730").Value = 2.05
$ws.Range("T730").Value = 1.75
$ws.Range("U730").Value = 1.8
$ws.Range("V730").Value = 2.05

# Row 731
$ws.Range("B731").Value = 7698680
$ws.Range("F731").Value = "CA Guemes"
$ws.Range("G731").Value = "Agropecuario"
$ws.Range("K731").Value = 2.2
$ws.Range("L731").Value = 3
$ws.Range("M731").Value = 3.2
$ws.Range("N731").Value = 2.6
$ws.Range("P731").Value = 3.1
$ws.Range("R731").Value = 1.85
$ws.Range("S731").Value = 2
$ws.Range("U731").Value = 1.875
$ws.Range("V731").Value = 1.975

# Row 732
$ws.Range("B732").Value = 7698551
$ws.Range("E732").Value = 45347.75
$ws.Range("F732").Value = "Racing de Cordoba"
$ws.Range("G732").Value = "Tristan Suarez"
$ws.Range("K732").Value = 1.8
$ws.Range("L732").Value = 3.2
$ws.Range("M732").Value = 4
$ws.Range("N732").Value = 1.95
$ws.Range("O732").Value = 3.3
$ws.Range("P732").Value = 4.2
$ws.Range("Q732").Value = -0.5
$ws.Range("R732").Value = 2
$ws.Range("S732").Value = 1.85
$ws.Range("T732").Value = 2
$ws.Range("U732").Value = 1.85
$ws.Range("V732").Value = 2

# Row 733
$ws.Range("B733").Value = 7698688
$ws.Range("E733").Value = 45347.77083333334
$ws.Range("F733").Value = "CA Chaco For Ever"
$ws.Range("G733").Value = "Deportivo Moron"
$ws.Range("K733").Value = 2.3
$ws.Range("L733").Value = 2.75
$ws.Range("M733").Value = 3.2
$ws.Range("N733").Value = 2.6
$ws.Range("O733").Value = 2.8
$ws.Range("P733").Value = 3.2
$ws.Range("Q733").Value = 0
$ws.Range("R733").Value = 1.8
$ws.Range("S733").Value = 2.05
$ws.Range("T733").Value = 1.75
$ws.Range("U733").Value = 1.875
$ws.Range("V733").Value = 1.975

# Row 734
$ws.Range("B734").Value = 7698691
$ws.Range("E734").Value = 45347.83333333334
$ws.Range("F734").Value = "Colon"
$ws.Range("G734").Value = "Patronato Parana"
$ws.Range("K734").Value = 1.727
$ws.Range("L734").Value = 3.5
$ws.Range("M734").Value = 4
$ws.Range("N734").Value = 1.6
$ws.Range("O734").Value = 4
$ws.Range("P734").Value = 5.5
$ws.Range("Q734").Value = -0.75
$ws.Range("T734").Value = 2.25
$ws.Range("U734").Value = 1.975
$ws.Range("V734").Value = 1.875

# Row 735
$ws.Range("B735").Value = 7698686
$ws.Range("F735").Value = "Atletico Rafaela"
$ws.Range("G735").Value = "Gimnasia Mendoza"
$ws.Range("K735").Value = 2.1
$ws.Range("L735").Value = 2.9
$ws.Range("M735").Value = 3.5
$ws.Range("N735").Value = 2
$ws.Range("O735").Value = 3.2
$ws.Range("P735").Value = 4.2
$ws.Range("Q735").Value = -0.5
$ws.Range("R735").Value = 2.05
$ws.Range("S735").Value = 1.8
$ws.Range("T735").Value = 1.75
$ws.Range("U735").Value = 1.775
$ws.Range("V735").Value = 2.1

# Row 736
$ws.Range("B736").Value = 7698679
$ws.Range("F736").Value = "CA Alvarado"
$ws.Range("G736").Value = "All Boys"
$ws.Range("N736").Value = 2.05
$ws.Range("O736").Value = 3
$ws.Range("P736").Value = 4.333
$ws.Range("U736").Value = 1.825
$ws.Range("V736").Value = 2.025

# Row 737
$ws.Range("B737").Value = 7698682
$ws.Range("E737").Value = 45347.90625
$ws.Range("F737").Value = "Quilmes"
$ws.Range("G737").Value = "Chacarita Juniors"
$ws.Range("K737").Value = 2.5
$ws.Range("M737").Value = 2.75
$ws.Range("N737").Value = 2.6
$ws.Range("P737").Value = 3
$ws.Range("Q737").Value = 0
$ws.Range("R737").Value = 1.8
$ws.Range("S737").Value = 2.05
$ws.Range("T737").Value = 2
$ws.Range("U737").Value = 2.05
$ws.Range("V737").Value = 1.8

# Row 738
$ws.Range("B738").Value = 7702207
$ws.Range("E738").Value = 45348.70833333334
$ws.Range("F738").Value = "CA San Miguel"
$ws.Range("G738").Value = "San Martin de Tucuman"
$ws.Range("K738").Value = 2.75
$ws.Range("L738").Value = 2.875
$ws.Range("M738").Value = 2.5
$ws.Range("N738").Value = 2.75
$ws.Range("O738").Value = 2.8
$ws.Range("R738").Value = 1.85
$ws.Range("S738").Value = 2
$ws.Range("T738").Value = 1.75
$ws.Range("U738").Value = 1.95
$ws.Range("V738").Value = 1.9

# Row 739
$ws.Range("B739").Value = 7845890
$ws.Range("E739").Value = 45349.88194444445
$ws.Range("F739").Value = "CA Atlanta"
$ws.Range("G739").Value = "CA Aldosivi"
$ws.Range("K739").Value = 2.1
$ws.Range("L739").Value = 3
$ws.Range("M739").Value = 3.3
$ws.Range("N739").Value = 1.95
$ws.Range("O739").Value = 3
$ws.Range("P739").Value = 3.75
$ws.Range("Q739").Value = -0.5
$ws.Range("R739").Value = 2.025
$ws.Range("S739").Value = 1.825
$ws.Range("T739").Value = 2
$ws.Range("U739").Value = 2.1
$ws.Range("V739").Value = 1.775

